$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "kWPqK338"
$ws.Range("B2").Value = 231004201
$ws.Range("C2").Value = "xydntam18"
$ws.Range("D2").Value = "SEt24!#j"
$ws.Range("F2").Value = "QmSGTMZp"
$ws.Range("G2").Value = "lajY"

# Row 3
$ws.Range("A3").Value = "jrpaB602"
$ws.Range("B3").Value = 231004200
$ws.Range("C3").Value = "mxgtktk14"
$ws.Range("D3").Value = "P!2d#7mN"
$ws.Range("F3").Value = "IMIHRjWV"
$ws.Range("G3").Value = "btGl"
